$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new mp3 lines (line_id, voice_id, name, text)
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "LcfcDJNUP1GQjkzn1xUU"
$ws.Range("C4").Value = "Emily"
$ws.Range("D4").Value = "I'm gonna make him an offer he can't refuse."

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "ODq5zmih8GrVes37Dizd"
$ws.Range("C5").Value = "Patrick"
$ws.Range("D5").Value = "Toto, I've a feeling we're not in Kansas anymore."

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "XB0fDUnXU5powFXDhCwa"
$ws.Range("C6").Value = "Charlotte"
$ws.Range("D6").Value = "There's no crying in baseball!"

# Column B needs to widen to fit the new (longer) voice_id values
$ws.Columns.Item(2).ColumnWidth = 20.41984375

# Update the active cell selection to match the target state
$ws.Range("E10").Select() | Out-Null
